$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: force a cell to plain text (so numeric-looking strings like
# "302.51" or "42.816.90" are not auto-coerced to numbers/dates by Excel), then drop
# the temporary "Text" number format back to the default "Normal" style so the saved
# cell carries no extra style index (matches the source cells, which are plain inline
# strings with no explicit style).

# Row swaps: coin identity + link + price + volume changed together (rows 11/12, 36/37, 50/51)
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.21%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0795"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0704"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.35%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.57%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.535.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.56%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.81%  "

# Price / volume refreshes for all other rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.816.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.302.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.670.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.296.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.796"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.800.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0900"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.52%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("E29").Value = "  -5.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.62%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.67%  "
$ws.Range("E41").Value = "  -5.03%  "
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("E43").Value = "  -5.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.973.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0284"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "
